# Feria Lagunitas de Puerto Montt - Betarraga
# A new weekly price observation is inserted as row 335, which pushes every
# existing row from 335 through 417 down by one (336 through 418), growing
# the used range from A1:R417 to A1:R418.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at 335 (shifts 335:417 -> 336:418, inherits the
# date number format from the neighbouring rows for column D).
$ws.Rows.Item(335).Insert()

# Populate the newly inserted row 335 with the new observation.
$ws.Range("A335").Value = 4
$ws.Range("B335").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C335").Value = "Los Lagos"
$ws.Range("D335").Value = 44943
$ws.Range("E335").Value = 10
$ws.Range("F335").Value = 100114014
$ws.Range("G335").Value = "Betarraga"
$ws.Range("H335").Value = "Sin especificar"
$ws.Range("I335").Value = "Primera"
$ws.Range("J335").Value = 1200
$ws.Range("K335").Value = 1000
$ws.Range("L335").Value = 1000
$ws.Range("M335").Value = 1000
$ws.Range("N335").Value = "$/paquete 5 unidades"
$ws.Range("O335").Value = "Provincia de Cautín"
$ws.Range("P335").Value = 200
$ws.Range("Q335").Value = 5
$ws.Range("R335").Value = "Hortaliza"
